$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "15.02.2025"
$ws.Range("B13").Value = "09:00"
$ws.Range("C13").Value = "68HS574"
$ws.Range("D13").Value = 5
$ws.Range("E13").Value = "11:30"
$ws.Range("F13").Value = 6
$ws.Range("G13").Value = "13:45"
$ws.Range("H13").Value = 14
$ws.Range("I13").Value = "14:30"
$ws.Range("J13").Value = 1
$ws.Range("K13").Value = "Melih Karaman"
$ws.Range("L13").Value = 9
$ws.Range("M13").Value = 10

$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "11.07.2025"
$ws.Range("B14").Value = "10:00"
$ws.Range("C14").Value = "45HD132"
$ws.Range("D14").Value = 564
$ws.Range("E14").Value = "12:00"
$ws.Range("F14").Value = 614
$ws.Range("G14").Value = "12:20"
$ws.Range("H14").Value = 617
$ws.Range("I14").Value = "14:30"
$ws.Range("J14").Value = 50
$ws.Range("K14").Value = "Ela karaman "
$ws.Range("L14").Value = 116
$ws.Range("M14").Value = 23
